$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.765134014560337
$ws.Range("J2").Value = -0.5513017252472258
$ws.Range("K2").Value = 0.1696252062987764

# Row 3
$ws.Range("I3").Value = -2.801104564858618
$ws.Range("J3").Value = -2.080177633312616
$ws.Range("K3").Value = -2.78912955792647

# Row 4
$ws.Range("H4").Value = -10.27025489571591
$ws.Range("I4").Value = -9.549327964169908
$ws.Range("J4").Value = -10.25827988878376

# Row 5
$ws.Range("G5").Value = 7.901006904421558
$ws.Range("H5").Value = 8.62193383596756
$ws.Range("I5").Value = 7.912981911353705
$ws.Range("K5").Value = 8.87489265094917

# Row 6
$ws.Range("F6").Value = -0.2386646425575917
$ws.Range("G6").Value = 0.4822622889884105
$ws.Range("H6").Value = -0.2266896356254442
$ws.Range("J6").Value = 0.73522110397002
$ws.Range("K6").Value = 1.048046298935466

# Row 7
$ws.Range("E7").Value = -2.386235114356984
$ws.Range("F7").Value = -1.665308182810982
$ws.Range("G7").Value = -2.374260107424837
$ws.Range("I7").Value = -1.412349367829372
$ws.Range("J7").Value = -1.099524172863926
$ws.Range("K7").Value = -1.977926754115506

# Row 8
$ws.Range("D8").Value = 1.053794872166346
$ws.Range("E8").Value = 1.774721803712348
$ws.Range("F8").Value = 1.065769879098493
$ws.Range("H8").Value = 2.027680618693958
$ws.Range("I8").Value = 2.340505813659403
$ws.Range("J8").Value = 1.462103232407824
$ws.Range("K8").Value = 0.7321912172293545

# Row 9
$ws.Range("C9").Value = 1.114969050580547
$ws.Range("D9").Value = 1.835895982126549
$ws.Range("E9").Value = 1.126944057512694
$ws.Range("G9").Value = 2.088854797108159
$ws.Range("H9").Value = 2.401679992073604
$ws.Range("I9").Value = 1.523277410822025
$ws.Range("J9").Value = 0.7933653956435556
$ws.Range("K9").Value = 2.455544492033183

# Row 10
$ws.Range("B10").Value = -0.9264868865757077
$ws.Range("C10").Value = -0.2055599550297054
$ws.Range("D10").Value = -0.9145118796435601
$ws.Range("F10").Value = 0.04739885995190407
$ws.Range("G10").Value = 0.36022405491735
$ws.Range("H10").Value = -0.5181785263342299
$ws.Range("I10").Value = -1.248090541512699
$ws.Range("J10").Value = 0.4140885548769285
$ws.Range("K10").Value = -0.2051307335183153

# Row 11
$ws.Range("B11").Value = 0.3770345820039356
$ws.Range("C11").Value = -0.3319173426099191
$ws.Range("E11").Value = 0.6299933969855451
$ws.Range("F11").Value = 0.942818591950991
$ws.Range("G11").Value = 0.06441601069941108
$ws.Range("H11").Value = -0.6654960044790579
$ws.Range("I11").Value = 0.9966830919105695
$ws.Range("J11").Value = 0.3774638035153257
$ws.Range("K11").Value = 0.6147675671350392

# Row 12
$ws.Range("B12").Value = -0.4275923834192769
$ws.Range("D12").Value = 0.5343183561761873
$ws.Range("E12").Value = 0.8471435511416332
$ws.Range("F12").Value = -0.03125903010994671
$ws.Range("G12").Value = -0.7611710452884157
$ws.Range("H12").Value = 0.9010080511012117
$ws.Range("I12").Value = 0.2817887627059679
$ws.Range("J12").Value = 0.5190925263256815

# Row 13
$ws.Range("C13").Value = 0.7947373931749101
$ws.Range("D13").Value = 1.107562588140356
$ws.Range("E13").Value = 0.2291600068887761
$ws.Range("F13").Value = -0.5007520082896928
$ws.Range("G13").Value = 1.161427088099934
$ws.Range("H13").Value = 0.5422077997046907
$ws.Range("I13").Value = 0.7795115633244043
$ws.Range("K13").Value = -0.1153642338804421

# Row 14
$ws.Range("B14").Value = -0.04071760298358112
$ws.Range("C14").Value = 0.2721075919818648
$ws.Range("D14").Value = -0.6062949892697151
$ws.Range("E14").Value = -1.336207004448184
$ws.Range("F14").Value = 0.3259720919414433
$ws.Range("G14").Value = -0.2932471964538005
$ws.Range("H14").Value = -0.05594343283408693
$ws.Range("J14").Value = -0.9508192300389333
$ws.Range("K14").Value = -0.3403303223714723

# Row 15
$ws.Range("B15").Value = 0.3721869518844864
$ws.Range("C15").Value = -0.5062156293670936
$ws.Range("D15").Value = -1.236127644545562
$ws.Range("E15").Value = 0.4260514518440648
$ws.Range("F15").Value = -0.193167836551179
$ws.Range("G15").Value = 0.04413592706853459
$ws.Range("I15").Value = -0.8507398701363118
$ws.Range("J15").Value = -0.2402509624688508
$ws.Range("K15").Value = -0.4017729932881683

# Row 16
$ws.Range("B16").Value = -0.1524291232873974
$ws.Range("C16").Value = -0.8823411384658664
$ws.Range("D16").Value = 0.779837957923761
$ws.Range("E16").Value = 0.1606186695285172
$ws.Range("F16").Value = 0.3979224331482308
$ws.Range("H16").Value = -0.4969533640566156
$ws.Range("I16").Value = 0.1135355436108454
$ws.Range("J16").Value = -0.04798648720847212

# Row 17
$ws.Range("B17").Value = -1.030518528898312
$ws.Range("C17").Value = 0.6316605674913157
$ws.Range("D17").Value = 0.0124412790960719
$ws.Range("E17").Value = 0.2497450427157855
$ws.Range("G17").Value = -0.6451307544890609
$ws.Range("H17").Value = -0.03464184682159993
$ws.Range("I17").Value = -0.1961638776409175

# Row 18
$ws.Range("B18").Value = 0.4742145784871607
$ws.Range("C18").Value = -0.1450047099080831
$ws.Range("D18").Value = 0.0922990537116305
$ws.Range("F18").Value = -0.8025767434932158
$ws.Range("G18").Value = -0.1920878358257549
$ws.Range("H18").Value = -0.3536098666450724

# Row 19
$ws.Range("B19").Value = 0.3556547466179877
$ws.Range("C19").Value = 0.5929585102377013
$ws.Range("E19").Value = -0.3019172869671451
$ws.Range("F19").Value = 0.3085716207003159
$ws.Range("G19").Value = 0.1470495898809984

# Row 20
$ws.Range("B20").Value = 0.3126006297022321
$ws.Range("D20").Value = -0.5822751675026142
$ws.Range("E20").Value = 0.02821374016484672
$ws.Range("F20").Value = -0.1333082906544708

# Row 21
$ws.Range("C21").Value = -0.4103003096576026
$ws.Range("D21").Value = 0.2001885980098584
$ws.Range("E21").Value = 0.03866656719054083

# Row 22
$ws.Range("B22").Value = -0.716162849403934
$ws.Range("C22").Value = -0.1056739417364731
$ws.Range("D22").Value = -0.2671959725557906

# Row 23
$ws.Range("B23").Value = 0.506656010950813
$ws.Range("C23").Value = 0.3451339801314955

# Row 24
$ws.Range("B24").Value = -0.343237405067616

